$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new data rows (256-269)
$ws.Cells.Item(256, 1).Value = 44330
$ws.Cells.Item(256, 2).Value = 5
$ws.Cells.Item(256, 3).Value = 7
$ws.Cells.Item(256, 4).Value = 85.07535245503161
$ws.Cells.Item(257, 1).Value = 44331
$ws.Cells.Item(257, 2).Value = 1
$ws.Cells.Item(257, 3).Value = 8
$ws.Cells.Item(257, 4).Value = 97.22897423432183
$ws.Cells.Item(258, 1).Value = 44332
$ws.Cells.Item(258, 2).Value = 2
$ws.Cells.Item(258, 3).Value = 10
$ws.Cells.Item(258, 4).Value = 121.5362177929023
$ws.Cells.Item(259, 1).Value = 44333
$ws.Cells.Item(259, 2).Value = 2
$ws.Cells.Item(259, 3).Value = 12
$ws.Cells.Item(259, 4).Value = 145.8434613514827
$ws.Cells.Item(260, 1).Value = 44334
$ws.Cells.Item(260, 2).Value = 4
$ws.Cells.Item(260, 3).Value = 15
$ws.Cells.Item(260, 4).Value = 182.3043266893534
$ws.Cells.Item(261, 1).Value = 44335
$ws.Cells.Item(261, 2).Value = 0
$ws.Cells.Item(261, 3).Value = 15
$ws.Cells.Item(261, 4).Value = 182.3043266893534
$ws.Cells.Item(262, 1).Value = 44336
$ws.Cells.Item(262, 2).Value = 0
$ws.Cells.Item(262, 3).Value = 14
$ws.Cells.Item(262, 4).Value = 170.1507049100632
$ws.Cells.Item(263, 1).Value = 44337
$ws.Cells.Item(263, 2).Value = 0
$ws.Cells.Item(263, 3).Value = 9
$ws.Cells.Item(263, 4).Value = 109.3825960136121
$ws.Cells.Item(264, 1).Value = 44338
$ws.Cells.Item(264, 2).Value = 0
$ws.Cells.Item(264, 3).Value = 8
$ws.Cells.Item(264, 4).Value = 97.22897423432183
$ws.Cells.Item(265, 1).Value = 44339
$ws.Cells.Item(265, 2).Value = 0
$ws.Cells.Item(265, 3).Value = 6
$ws.Cells.Item(265, 4).Value = 72.92173067574137
$ws.Cells.Item(266, 1).Value = 44340
$ws.Cells.Item(266, 2).Value = 0
$ws.Cells.Item(266, 3).Value = 4
$ws.Cells.Item(266, 4).Value = 48.61448711716091
$ws.Cells.Item(267, 1).Value = 44341
$ws.Cells.Item(267, 2).Value = 0
$ws.Cells.Item(267, 3).Value = 0
$ws.Cells.Item(267, 4).Value = 0
$ws.Cells.Item(268, 1).Value = 44342
$ws.Cells.Item(268, 2).Value = 0
$ws.Cells.Item(268, 3).Value = 0
$ws.Cells.Item(268, 4).Value = 0
$ws.Cells.Item(269, 1).Value = 44343
$ws.Cells.Item(269, 2).Value = 0
$ws.Cells.Item(269, 3).Value = 0
$ws.Cells.Item(269, 4).Value = 0

# Copy formatting from the last pre-existing row's date cell (A255) down to the
# new rows' date cells (A256:A269) so they pick up style index 2 (date format),
# matching the existing pattern in column A. Columns B-D keep the default style.
$ws.Range("A255").Copy()
$ws.Range("A256:A269").PasteSpecial(-4122)

$excel.CutCopyMode = 0
